$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Both are fine.', 'どちらでもいいです。')
    ,@('Same thing.', '同じです。|おなじです。')
    ,@('More or less the same.', 'だいたい同じです。|だいたいおなじです。')
    ,@('A little different.', 'ちょっと違います。|ちょっとちがいます。')
    ,@('Can''t use it.', '使えません。|つかえません。')
    ,@('No good.', 'だめです。')
    ,@('Raise your hand.', '手をあげてください。|てをあげてください。')
    ,@('Read it before coming to class.', '読んできてください。|よんできてください。')
    ,@('Hand in the homework.', '宿題を出してください。|しゅくだいをだしてください。')
    ,@('Open the book to page 10.', '10ページを開いてください。|10ページをひらいてください。')
    ,@('Close the textbook.', '教科書を閉じてください。|きょうかしょをとじてください。')
    ,@('Ask the person sitting next to you.', 'となりの人に聞いてください。|となりのひとにきいてください。')
    ,@('Please stop.', 'やめてください。')
    ,@('That''s it for today.', '今日はこれで終わります。|きょうはこれでおわります。')
    ,@('deadline', 'しめきり')
    ,@('exercise', '練習|れんしゅう')
    ,@('meaning', '意味|いみ')
    ,@('pronunciation', '発音|はつおん')
    ,@('grammar', '文法|ぶんぽう')
    ,@('question', '質問|しつもん')
    ,@('answer', '答|こたえ')
    ,@('example', '例|れい')
    ,@('( ) (parenthesis)', 'かっこ')
    ,@('〇 (correct)', 'まる')
    ,@('✕ (wrong)', 'ばつ')
    ,@('colloquial expression', 'くだけた言い方|くだけたいいかた')
    ,@('bookish expression', 'かたい言い方|かたいいいかた')
    ,@('polite expression', 'ていねいな言い方|ていねいないいかた')
    ,@('dialect', '方言|ほうげん')
    ,@('common language', '共通語|きょうつうご')
    ,@('for example', 'たとえば')
    ,@('anything else', 'ほかに')
    ,@('number', '～番|～ばん')
    ,@('line number', '～行目|～ぎょうめ')
    ,@('two people each', '二人ずつ|ふたりずつ')
)

$startRow = 65
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
